# Baseline LH_REVIEWS Folder from dev branch
# Adds a new "Version History" entry (row 5) documenting the v2.0 review
# pass performed by Hala Eldaly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Version History")

$ws.Range("A5").Value = "v2.0"
$ws.Range("B5").Value = "Hala Eldaly"
$ws.Range("C5").Value = "Review DELETPOST And No Comments "
$ws.Range("D5").Value = Get-Date -Year 2025 -Month 5 -Day 13 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Keep the active sheet / selection consistent with the authored workbook.
$ws.Activate()
$ws.Range("A1:XFD1048576").Select()
